$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row: add columns O1:R1 = 13,14,15,16, matching N1 style (bold/center/border) ---
$hdr = New-Object 'object[,]' 1,4
$hdr[0,0] = 13
$hdr[0,1] = 14
$hdr[0,2] = 15
$hdr[0,3] = 16
$ws.Range("O1:R1").Value = $hdr
$ws.Range("N1").Copy()
$ws.Range("O1:R1").PasteSpecial(-4122)

# --- Block 1: B2:G25 ---
$csv1 = @"
1.05,1.012443665433292,1.033256565731016,1.02710283525598,1.037751313904113,1.0
1.05,1.016295163003192,1.035862927170052,1.030080779886116,1.040573060101994,1.0
1.05,1.018743663555657,1.03752139906953,1.031979759153774,1.042372647957879,1.0
1.05,1.01976583381299,1.038215655663923,1.032774334389188,1.043125693966674,1.0
1.05,1.019940481655931,1.0383362038533,1.032910719943221,1.043254968921108,1.0
1.05,1.018767039119968,1.037542359715049,1.031999392201719,1.042391274739579,1.0
1.05,1.013766521562679,1.034157580508013,1.028126298866642,1.038720982835876,1.0
1.049999999999999,1.004590421631043,1.027952862691021,1.021064474035381,1.032031030900262,1.0
1.05,0.9982774490548367,1.023712026609233,1.016264112975261,1.027516647306324,1.0
1.05,0.9959521143004051,1.022266473934344,1.014688881335861,1.026287480476826,1.0
1.05,0.9952444265008984,1.021870408153502,1.014288956878493,1.02609746311887,1.0
1.05,0.9957504751830173,1.022273335037862,1.014781724611238,1.026715696734136,1.0
1.05,0.9966549785875943,1.022918313266063,1.015530988688195,1.027514829490795,1.0
1.05,0.9971449096444982,1.023256493439041,1.015915778706239,1.027897747771318,1.0
1.05,0.999727181101508,1.024985332056787,1.017861209220292,1.029718204654549,1.0
1.05,1.001235035974094,1.025973945776008,1.01896187601687,1.030692605517853,1.0
1.05,1.001947413476282,1.026399091654215,1.019415060704292,1.030992881829269,1.0
1.05,1.001945934004852,1.026323827318407,1.019285517983431,1.030685522580288,1.0
1.05,0.9999386556879678,1.024839324081886,1.017529883348571,1.028708101377491,1.0
1.05,0.9951067880232575,1.021574518043839,1.013832130520935,1.025182061669854,1.0
1.05,0.9920273578751614,1.019499955111555,1.011493535244187,1.022968755752015,1.0
1.05,0.9936570046801472,1.020593355707911,1.012728825779209,1.024138083451543,1.0
1.05,0.9999585862844801,1.024835327555667,1.017524683303721,1.028678473451442,1.0
1.05,1.007021362601134,1.029602132495484,1.022931897665342,1.033799604965994,1.0
"@
$lines = $csv1 -split "`n" | Where-Object { $_.Trim().Length -gt 0 }
$nrows = $lines.Count
$ncols = ($lines[0].Trim() -split ",").Count
$b1 = New-Object 'object[,]' $nrows,$ncols
for ($i = 0; $i -lt $nrows; $i++) {
    $parts = $lines[$i].Trim() -split ","
    for ($j = 0; $j -lt $ncols; $j++) {
        $b1[$i,$j] = $parts[$j] -as [double]
    }
}
$ws.Range("B2:G25").Value = $b1

# --- Block 2: I2:N25 ---
$csv2 = @"
1.053397357240914,1.034258453584263,1.044271239482289,1.038197068704896,1.048708621095344,1.005712725503999
1.054231618556763,1.036351251793774,1.046056611491831,1.040342711859519,1.050711842858733,1.005712725503983
1.05474823344294,1.037679760418715,1.047187120357766,1.041706949491678,1.051985216776259,1.005712725503983
1.054961756222442,1.038235163291501,1.047660131487667,1.042277561625228,1.052517798670675,1.005712725503983
1.054999506204118,1.038331593123236,1.047743528156205,1.042376356706941,1.052610068153125,1.005712725503983
1.054757038716128,1.037696647532572,1.047204993378899,1.041723486201667,1.052000801655056,1.005712725503983
1.053691404882493,1.034982785801368,1.044893920420342,1.03893829676368,1.049400804918378,1.005712725503983
1.051623278790479,1.029982115175805,1.040608138213175,1.033824795488659,1.044624956464611,1.005712725503983
1.050146130972614,1.026559577176331,1.037665956477846,1.030345794465968,1.041406239443575,1.005712725503983
1.049693426155365,1.025494265043817,1.036780348738243,1.029339003339359,1.040730115827387,1.005712725503983
1.049594835885812,1.025251423353936,1.036590267234168,1.02914742932185,1.040741193241993,1.005712725503983
1.049776853850158,1.025642494471923,1.036943257349244,1.029588003518182,1.041305982607775,1.005712725503983
1.050031548363608,1.026203433505226,1.037437009991387,1.030182431621828,1.041952178007195,1.005712725503983
1.050158324173412,1.026485136912159,1.037683192390314,1.030473365236317,1.042242940801339,1.005712725503983
1.050771824884016,1.027874807496442,1.038882263519461,1.031879539181466,1.043535682202674,1.005712725503983
1.051105268442126,1.028645836163369,1.039541385073802,1.03264540731582,1.044183064420459,1.005712725503983
1.05122198950283,1.028937051217506,1.039777186853889,1.032906873966168,1.04429725496655,1.005712725503983
1.051150066391273,1.028801428738408,1.03964094446132,1.032716655082515,1.043932981934789,1.005712725503983
1.050548860934962,1.027471942260317,1.038459026798447,1.03127133877904,1.042264262272799,1.005712725503983
1.049373038196806,1.024810526646969,1.036158366644916,1.028556103449401,1.03970155768443,1.005712725503983
1.048613490370378,1.023121510868567,1.034694765243035,1.026840690499261,1.038098528216405,1.005712725503983
1.049013794426864,1.024011854729415,1.035463973741116,1.027745198632784,1.038943958604129,1.005712725503983
1.050538682516767,1.027458230487882,1.038439906043295,1.031250879914058,1.042220009202578,1.005712725503983
1.052188911029614,1.03131471854805,1.041757658197859,1.035184357803129,1.045894942886257,1.005712725503983
"@
$lines = $csv2 -split "`n" | Where-Object { $_.Trim().Length -gt 0 }
$nrows = $lines.Count
$ncols = ($lines[0].Trim() -split ",").Count
$b2 = New-Object 'object[,]' $nrows,$ncols
for ($i = 0; $i -lt $nrows; $i++) {
    $parts = $lines[$i].Trim() -split ","
    for ($j = 0; $j -lt $ncols; $j++) {
        $b2[$i,$j] = $parts[$j] -as [double]
    }
}
$ws.Range("I2:N25").Value = $b2

# --- Block 3: O2:R25 (new columns) ---
$csv3 = @"
1.03,1.047121966925635,1.02,1.042373881393532
1.03,1.048707366549118,1.02,1.04363363127527
1.03,1.049715144638864,1.02,1.044433880136991
1.03,1.050136642081041,1.02,1.044775459785955
1.03,1.050209666226928,1.02,1.044843115802429
1.03,1.049727478870314,1.02,1.044466503332012
1.03,1.047669778836184,1.02,1.042836946519865
1.03,1.043890042408065,1.02,1.039803479440762
1.03,1.041393908187833,1.02,1.037739935905935
1.03,1.041293630055633,1.02,1.037146520661917
1.03,1.041627819559078,1.02,1.037012127617748
1.03,1.042350309786162,1.02,1.03725921646957
1.03,1.043034007419545,1.02,1.037609734120478
1.03,1.043301275433228,1.02,1.037789644162284
1.03,1.044284439042242,1.02,1.038640578851898
1.03,1.044667475485918,1.02,1.039109176097956
1.03,1.044520831528573,1.02,1.039264302664003
1.03,1.0439085059227,1.02,1.039174345902155
1.03,1.042062428351283,1.02,1.038342599171675
1.03,1.039993510139087,1.02,1.036719178906716
1.03,1.038724817782537,1.02,1.035670883990347
1.03,1.039393920590844,1.02,1.036205131014008
1.03,1.041986694306046,1.02,1.038301765400091
1.03,1.044895146121197,1.02,1.040644716243662
"@
$lines = $csv3 -split "`n" | Where-Object { $_.Trim().Length -gt 0 }
$nrows = $lines.Count
$ncols = ($lines[0].Trim() -split ",").Count
$b3 = New-Object 'object[,]' $nrows,$ncols
for ($i = 0; $i -lt $nrows; $i++) {
    $parts = $lines[$i].Trim() -split ","
    for ($j = 0; $j -lt $ncols; $j++) {
        $b3[$i,$j] = $parts[$j] -as [double]
    }
}
$ws.Range("O2:R25").Value = $b3
